$wb = $excel.ActiveWorkbook

# --- Personnel sheet: correct Kevin Cahill's role from "technician" to "creator" ---
$ws = $wb.Worksheets.Item("Personnel")
$ws.Range("G6").Value = "creator"

# --- Update the last active selection on the Personnel sheet (G7) ---
$ws.Range("G7").Select()
